# Auto-generated edit script: appends the "first 6 rows of Baltimore"
# source-code / output paragraphs to the end of the document, matching
# the target diff.

function Get-EndRange($doc) {
    $e = $doc.Content.End
    return $doc.Range($e, $e)
}

function Escape-Xml($text) {
    $t = $text
    $t = $t.Replace("&", "&amp;")
    $t = $t.Replace("<", "&lt;")
    $t = $t.Replace(">", "&gt;")
    return $t
}

# Appends one new paragraph (style "SourceCode") built from an ordered
# array of run descriptors:
#   @{Kind="text"; Text="..."; Style="CharStyleId"}
#   @{Kind="break"}                                   (manual line break)
# The paragraph/run skeleton is inserted as plain WordprocessingML via
# Range.InsertXML (character-style rStyle elements are not honored by
# InsertXML in this host), then each run's character style is applied
# afterwards via Range.Style so the saved OOXML carries the correct
# <w:rStyle> elements.
function Add-SourceParagraph($doc, $runs) {
    $xml = '<w:p><w:pPr><w:pStyle w:val="SourceCode"/></w:pPr>'
    foreach ($run in $runs) {
        if ($run.Kind -eq "break") {
            $xml = $xml + '<w:r><w:br/></w:r>'
        } else {
            $escaped = Escape-Xml $run.Text
            $xml = $xml + '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
        }
    }
    $xml = $xml + '</w:p>'

    $insertRange = Get-EndRange $doc
    $insertRange.InsertXML($xml)

    $paraCount = $doc.Paragraphs.Count
    $para = $doc.Paragraphs.Item($paraCount)
    $pos = $para.Range.Start

    foreach ($run in $runs) {
        if ($run.Kind -eq "break") {
            $pos = $pos + 1
        } else {
            $len = $run.Text.Length
            if ($run.Style) {
                $styledRange = $doc.Range($pos, $pos + $len)
                $styledRange.Style = $run.Style
            }
            $pos = $pos + $len
        }
    }
}

$d = $word.ActiveDocument

Add-SourceParagraph $d @(
    @{Kind="text"; Text='Arrest_of_Grey'; Style='NormalTok'},
    @{Kind="text"; Text='<-'; Style='OtherTok'},
    @{Kind="text"; Text=' homicide '; Style='NormalTok'},
    @{Kind="text"; Text='%>%'; Style='SpecialCharTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="break"},
    @{Kind="text"; Text='  '; Style='NormalTok'},
    @{Kind="text"; Text='filter'; Style='FunctionTok'},
    @{Kind="text"; Text='(victim_last '; Style='NormalTok'},
    @{Kind="text"; Text='=='; Style='SpecialCharTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="text"; Text='"GREY"'; Style='StringTok'},
    @{Kind="text"; Text=', city '; Style='NormalTok'},
    @{Kind="text"; Text='=='; Style='SpecialCharTok'},
    @{Kind="text"; Text='"Baltimore"'; Style='StringTok'},
    @{Kind="text"; Text=')'; Style='NormalTok'}
)

Add-SourceParagraph $d @(
    @{Kind="text"; Text='library'; Style='FunctionTok'},
    @{Kind="text"; Text='(lubridate)'; Style='NormalTok'}
)

Add-SourceParagraph $d @(
    @{Kind="text"; Text='## '; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='## Attaching package: ''lubridate'''; Style='VerbatimChar'}
)

Add-SourceParagraph $d @(
    @{Kind="text"; Text='## The following objects are masked from ''package:base'':'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='## '; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='##     date, intersect, setdiff, union'; Style='VerbatimChar'}
)

Add-SourceParagraph $d @(
    @{Kind="text"; Text='Baltimore '; Style='NormalTok'},
    @{Kind="text"; Text='<-'; Style='OtherTok'},
    @{Kind="text"; Text=' homicide '; Style='NormalTok'},
    @{Kind="text"; Text='%>%'; Style='SpecialCharTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="break"},
    @{Kind="text"; Text='  '; Style='NormalTok'},
    @{Kind="text"; Text='mutate'; Style='FunctionTok'},
    @{Kind="text"; Text='('; Style='NormalTok'},
    @{Kind="text"; Text='reported_date ='; Style='AttributeTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="text"; Text='ymd'; Style='FunctionTok'},
    @{Kind="text"; Text='(reported_date)) '; Style='NormalTok'},
    @{Kind="text"; Text='%>%'; Style='SpecialCharTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="break"},
    @{Kind="text"; Text='  '; Style='NormalTok'},
    @{Kind="text"; Text='filter'; Style='FunctionTok'},
    @{Kind="text"; Text='(city_name '; Style='NormalTok'},
    @{Kind="text"; Text='=='; Style='SpecialCharTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="text"; Text='''Baltimore,MD'''; Style='StringTok'},
    @{Kind="text"; Text=') '; Style='NormalTok'},
    @{Kind="text"; Text='%>%'; Style='SpecialCharTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="break"},
    @{Kind="text"; Text='  '; Style='NormalTok'},
    @{Kind="text"; Text='group_by'; Style='FunctionTok'},
    @{Kind="text"; Text='('; Style='NormalTok'},
    @{Kind="text"; Text='date ='; Style='AttributeTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="text"; Text='floor_date'; Style='FunctionTok'},
    @{Kind="text"; Text='(reported_date, '; Style='NormalTok'},
    @{Kind="text"; Text='''month'''; Style='StringTok'},
    @{Kind="text"; Text=')) '; Style='NormalTok'},
    @{Kind="text"; Text='%>%'; Style='SpecialCharTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="break"},
    @{Kind="text"; Text='  '; Style='NormalTok'},
    @{Kind="text"; Text='summarize'; Style='FunctionTok'},
    @{Kind="text"; Text='('; Style='NormalTok'},
    @{Kind="text"; Text='homicide ='; Style='AttributeTok'},
    @{Kind="text"; Text=' '; Style='NormalTok'},
    @{Kind="text"; Text='n'; Style='FunctionTok'},
    @{Kind="text"; Text='())'; Style='NormalTok'}
)

Add-SourceParagraph $d @(
    @{Kind="text"; Text='## Warning: 2 failed to parse.'; Style='VerbatimChar'}
)

Add-SourceParagraph $d @(
    @{Kind="text"; Text='head'; Style='FunctionTok'},
    @{Kind="text"; Text='(Baltimore)'; Style='NormalTok'}
)

Add-SourceParagraph $d @(
    @{Kind="text"; Text='## # A tibble: 6 × 2'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='##   date       homicide'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='##   <date>        <int>'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='## 1 2007-01-01       28'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='## 2 2007-02-01       17'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='## 3 2007-03-01       26'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='## 4 2007-04-01       19'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='## 5 2007-05-01       32'; Style='VerbatimChar'},
    @{Kind="break"},
    @{Kind="text"; Text='## 6 2007-06-01       31'; Style='VerbatimChar'}
)


Write-Host "Paragraph count after edit: $($d.Paragraphs.Count)"
